$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text formatting (these values look
# numeric/percentage-like, e.g. "563.77", "0.520", "  +3.70%  ") by forcing a
# text number format before assigning the value, matching the source data
# which stores all of these as plain text (inline strings).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.656.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.001.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.77"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.62%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.990.88"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.134"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.98%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.97"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.499.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.002.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.658.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "437.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.17%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.44"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +10.84%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0782"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +15.24%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.991"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.58"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.40"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0356"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.767.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.31%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.56"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +20.12%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.57%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.67"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.59%  "
